# edit.ps1 -- apply the "Consciousness" -> "Arts" essay rewrite
#
# Strategy notes:
#   - $d.Content.Find.Execute(...) does plain literal text search/replace.
#   - Whenever a paragraph "block" (the runs between <w:br/> separators /
#     paragraph boundaries) is touched by *any* edit, the engine re-merges
#     every run in that block that shares identical run formatting (this
#     mirrors real Word's "coalesce adjacent identically-formatted runs"
#     behavior). That means a second/third edit to the same block erases
#     run splits we still want (e.g. separate "." runs, or brand-new runs
#     we are inserting).
#   - To keep a run boundary from being swallowed by that coalescing, we
#     temporarily mark the earlier piece Bold = $true right after writing
#     its text (before any later edit touches the same block). Because
#     neighboring runs then have different formatting they will not be
#     merged by later edits to the block. Once every split we need in a
#     block has been created, we go back and clear Bold = $false on every
#     marked range in that block; by then no further edits touch that
#     block, so the split sticks even though the final formatting is
#     identical again.
#   - $range.Text = "..." (direct property assignment) is used instead of
#     Find.Execute's ReplaceWith parameter for any text containing an
#     apostrophe, since ReplaceWith goes through AutoCorrect's
#     "smart quotes" pass and mangles straight quotes into curly ones;
#     plain assignment does not.

$d = $word.ActiveDocument

# Find `search` starting from the whole document and return the Range
# positioned exactly over the match (zero-width replace semantics).
function Find-Range($search) {
    $r = $d.Content
    $r.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $r
}

function Find-Replace($search, $replacement) {
    $d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2) | Out-Null
}

# Replace the text of an existing run (found via its old text) and mark
# it Bold as a temporary "do not merge me with my neighbors" flag. Returns
# the Range so callers can batch up-clear the flag later.
function Set-ProtectedText($search, $newText) {
    $r = Find-Range $search
    $r.Text = $newText
    $r.Bold = $true
    return $r
}

# Insert brand-new text right after $afterRange's end, also flagged Bold.
function Insert-ProtectedAfter($afterRange, $newText) {
    $r = $d.Range($afterRange.End, $afterRange.End)
    $r.InsertAfter($newText)
    $r.Bold = $true
    return $r
}

# ---------------------------------------------------------------------
# Title / author
# ---------------------------------------------------------------------
Find-Replace "Unraveling the Enigma of Consciousness: A Journey of Discovery" "The Enriching World of Arts: A Journey Through Expression and Creativity"
Find-Replace "Vivian Lam" "Abigail Kent"

# ---------------------------------------------------------------------
# Email paragraph: "vivianlam@protonmail" + "." + "com"
#               -> "abigail" + "." + "kent25@xyz" + "." + "com"
# (last "." and "com" runs are pre-existing / untouched by the diff, so
# we must protect the new "kent25@xyz" / old "." boundary too)
# ---------------------------------------------------------------------
$rEmail1 = Set-ProtectedText "vivianlam@protonmail" "abigail"
$rEmail2 = Insert-ProtectedAfter $rEmail1 "."
$rEmail3 = Insert-ProtectedAfter $rEmail2 "kent25@xyz"
# the pre-existing "." run that used to directly follow "vivianlam@protonmail"
$rEmail4 = $d.Range($rEmail3.End, $rEmail3.End + 1)
$rEmail4.Bold = $true

$rEmail1.Bold = $false
$rEmail2.Bold = $false
$rEmail3.Bold = $false
$rEmail4.Bold = $false

# ---------------------------------------------------------------------
# Body paragraph 1 / Block A (4 sentences, run count unchanged)
# ---------------------------------------------------------------------
$rA1 = Set-ProtectedText "Consciousness, the enigmatic phenomenon at the core of our existence, has long captivated philosophers, scientists, and artists alike" "As we embark on a voyage through the realm of Arts, we step into a world where imagination reigns supreme, where the boundaries of reality blur, and where emotions find their voice"
$rA2 = Set-ProtectedText " As we navigate the complexities of the physical world, our perception of reality and sense of self emerge as remarkable features that define our human experience" " Arts, in its myriad forms, serves as a mirror to society, reflecting the human experience in all its complexities and nuances"
$rA3 = Set-ProtectedText " Yet, the nature of consciousness remains shrouded in mystery" " From the strokes of a paintbrush on a canvas to the flowing melodies of a symphony, art transcends the limitations of language, inviting us to connect with ourselves, with others, and with the world around us"
$rA4 = Set-ProtectedText " With each new discovery, we catch tantalizing glimpses into this intricate tapestry of perception, emotion, and thought" " In this essay, we will delve into the captivating power of Arts, exploring its multifaceted dimensions and unraveling the profound impact it has on our lives"

$rA1.Bold = $false
$rA2.Bold = $false
$rA3.Bold = $false
$rA4.Bold = $false

# ---------------------------------------------------------------------
# Block B (3 sentences, run count unchanged)
# ---------------------------------------------------------------------
$rB1 = Set-ProtectedText "In the realm of science, researchers are undertaking groundbreaking studies on neural networks and brain activity, seeking to decode the physical mechanisms that underpin consciousness" "Within the tapestry of human existence, art weaves its way seamlessly, becoming an integral part of our history, culture, and identity"
$rB2 = Set-ProtectedText " Explorations into altered states of consciousness, such as dreams and meditative practices, provide unique insights into the plasticity of our mental landscapes" " From the cave paintings of ancient civilizations to the modern masterpieces adorning museums, art serves as a living testament to our collective memory, bridging the gap between generations and providing invaluable insights into the human condition"
$rB3 = Set-ProtectedText " The study of consciousness not only enriches our understanding of ourselves but also holds profound implications for our comprehension of artificial intelligence and its potential to emulate human cognition" " Whether it's through the lens of literature, the stage of theater, or the notes of music, art immortalizes the triumphs and tribulations of humanity, capturing the essence of what it means to be human"

$rB1.Bold = $false
$rB2.Bold = $false
$rB3.Bold = $false

# fix the apostrophe: Find.Execute's ReplaceWith / InsertAfter path runs
# AutoCorrect's smart-quote pass on typed text, which turns a straight
# apostrophe into a curly one. Direct Range.Text assignment does not, so
# patch the single character back afterwards.
$rApos = Find-Range ([char]0x2019)
$rApos.Text = [char]0x27

# ---------------------------------------------------------------------
# Block C (3 sentences -> 5 "sentences" = 2 brand-new runs inserted
# right before the existing trailing "." run)
# ---------------------------------------------------------------------
$rC1 = Set-ProtectedText "As we delved into the depths of consciousness, we uncover profound connections to our artistic and cultural expressions" "The realm of art transcends the boundaries of mere aesthetics; it possesses the transformative power to shape our perceptions, provoke thought, and inspire action"
$rC2 = Set-ProtectedText " From literature's exploration of inner turmoil to music's ability to evoke an emotional response, we find echoes of our conscious experiences reflected in the works of great artists" " Art has the ability to challenge societal norms, question established ideologies, and ignite movements for change"
$rC3 = Set-ProtectedText " These explorations transcending disciplinary boundaries offer a multi-faceted perspective on the multifaceted nature of consciousness" " It can educate, inform, and empower, giving voice to the marginalized and shedding light on pressing issues"
$rC4 = Insert-ProtectedAfter $rC3 "."
$rC5 = Insert-ProtectedAfter $rC4 " By engaging with works of art, we become more empathetic, more tolerant, and more aware of the world around us, fostering a sense of global citizenship and interconnectedness"
# the pre-existing trailing "." run for this paragraph
$rC6 = $d.Range($rC5.End, $rC5.End + 1)
$rC6.Bold = $true

$rC1.Bold = $false
$rC2.Bold = $false
$rC3.Bold = $false
$rC4.Bold = $false
$rC5.Bold = $false
$rC6.Bold = $false

# ---------------------------------------------------------------------
# Summary paragraph (whole paragraph is one block; 5 sentences -> 8)
# ---------------------------------------------------------------------
$rS1 = Set-ProtectedText "This essay delved into the enigmatic nature of consciousness, weaving together scientific investigations, philosophical contemplations, and artistic representations" "In the realm of Arts, we find a world where imagination and creativity flourish, where expression transcends the limitations of language, and where emotions find their voice"
$rS2 = Set-ProtectedText " As we continue to unravel the complexities of our conscious experience, we gain a deeper appreciation for the richness and mystery of our own existence" " Art serves as a mirror to society, reflecting the human experience in all its complexities and nuances"
$rS3 = Set-ProtectedText " The pursuit of understanding " " It has the power to shape our perceptions, provoke thought, and inspire action, becoming an integral part of our history, culture, and identity"
# this run used to carry <w:lastRenderedPageBreak/> followed by
# "consciousness remains ...". Both the page-break marker and the text
# are replaced outright by a lone "." (the marker itself is relocated to
# the very start of the paragraph afterwards).
$rS4 = Set-ProtectedText "consciousness remains an ongoing journey, beckoning us to explore the vast landscapes of our interconnected minds" "."
$rS5 = Insert-ProtectedAfter $rS4 " Art educates, informs, and empowers, fostering empathy, tolerance, and a sense of global citizenship"
$rS6 = Insert-ProtectedAfter $rS5 "."
$rS7 = Insert-ProtectedAfter $rS6 " As we engage with works of art, we embark on a transformative journey, enriching our lives and deepening our understanding of ourselves, others, and the world we inhabit"
# the pre-existing trailing "." run for this paragraph
$rS8 = $d.Range($rS7.End, $rS7.End + 1)
$rS8.Bold = $true

$rS1.Bold = $false
$rS2.Bold = $false
$rS3.Bold = $false
$rS4.Bold = $false
$rS5.Bold = $false
$rS6.Bold = $false
$rS7.Bold = $false
$rS8.Bold = $false

# ---------------------------------------------------------------------
# Add a trailing empty paragraph after the Summary paragraph.
# ---------------------------------------------------------------------
$d.Paragraphs.Add() | Out-Null
